$wb = $excel.ActiveWorkbook

$rowsData = @(
    @{ Row=2; B='2024-02-12'; C='赣州·宅舞联萌·随舞动漫派对（免费活动)'; D='金岭东大道新都汇西侧约100米 万达广场'; E='2024.02.12 14:00-02.13 19:00'; F=56; G='不可售'; H='https://show.bilibili.com/platform/detail.html?id=81540'; I='//i0.hdslb.com/bfs/openplatform/202401/5gLDYtbv1706608938962.jpeg' },
    @{ Row=3; B='2024-02-14'; C=' 赣州·十万伏特·2024次元交流会（冬）'; D='平安大道 麋鹿LiveHose'; E='2024.02.14 09:30-02.15 17:30'; F=106; G=35; H='https://show.bilibili.com/platform/detail.html?id=81248'; I='//i0.hdslb.com/bfs/openplatform/202401/mKDiDPv31705921109896.jpeg' },
    @{ Row=4; B='2024-02-14'; C='南昌·原X穹X崩only'; D='龙蟠街666号融创茂1层 融创茂'; E='2024.02.14 10:00-02.15 17:00'; F=141; G=60; H='https://show.bilibili.com/platform/detail.html?id=80784'; I='//i2.hdslb.com/bfs/openplatform/202401/iNAvP52t1705039345817.jpeg' },
    @{ Row=5; B='2024-02-14'; C='南昌·龙年动漫展'; D='南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'; E='2024.02.14 10:00-02.15 18:00'; F=1363; G=55; H='https://show.bilibili.com/platform/detail.html?id=80525'; I='//i2.hdslb.com/bfs/openplatform/202401/ezt7koZo1704444854691.jpeg' },
    @{ Row=6; B='2024-02-14'; C='吉安·COMIC LIFE 次元假日03'; D='东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'; E='2024.02.14 09:30-02.14 18:00'; F=1602; G=45; H='https://show.bilibili.com/platform/detail.html?id=80305'; I='//i2.hdslb.com/bfs/openplatform/202312/gUyU7wgj1703754978855.jpeg' },
    @{ Row=7; B='2024-02-14'; C='宜春·第三十届静卿国风动漫新春盛典'; D='宜阳大道19号(交通银行旁) 宜春安缦文华酒店'; E='2024.02.14 09:00-02.14 17:00'; F=350; G=50; H='https://show.bilibili.com/platform/detail.html?id=80802'; I='//i1.hdslb.com/bfs/openplatform/202401/Pq8w7EsS1705048754533.jpeg' },
    @{ Row=8; B='2024-02-14'; C='江西·第二十一届九江ACJJ国际动漫展'; D='体育路九江市体育中心-体育馆 九江市体育中心'; E='2024.02.14 09:00-02.15 17:00'; F=462; G=50; H='https://show.bilibili.com/platform/detail.html?id=81015'; I='//i2.hdslb.com/bfs/openplatform/202401/HD1sIIY21705557926335.jpeg' },
    @{ Row=9; B='2024-02-14'; C='赣州·国乙ONLY（取消）'; D='八一四大道18号 纽豪花园酒店'; E='2024.02.14 10:00-02.14 17:00'; F=182; G='不可售'; H='https://show.bilibili.com/platform/detail.html?id=79487'; I='//i0.hdslb.com/bfs/openplatform/202312/Rhqex68Z1701935089796.jpeg' },
    @{ Row=10; B='2024-02-15'; C='萍乡·au7新年国漫展'; D='建设西路钻石公寓西南侧60米 智博篮球馆'; E='2024.02.15 10:00-02.15 17:00'; F=189; G=45; H='https://show.bilibili.com/platform/detail.html?id=80790'; I='//i2.hdslb.com/bfs/openplatform/202401/LiamxFS81705481738724.jpeg' },
    @{ Row=11; B='2024-02-15'; C='赣州·明日方舟ONLY大炎新岁同好交流茶话会'; D='南门口地一大道下沉广场 漫库书店'; E='2024.02.15 11:00-02.15 18:00'; F=167; G=48; H='https://show.bilibili.com/platform/detail.html?id=78689'; I='//i1.hdslb.com/bfs/openplatform/202311/T1Y8Iju31700621742031.png' },
    @{ Row=12; B='2024-02-15'; C='鹰潭·ADO7新年 原·星·蔚蓝 主题展'; D='南站路锦都金源酒店18楼 锦都金源酒店'; E='2024.02.15 10:00-02.15 17:00'; F=75; G=55; H='https://show.bilibili.com/platform/detail.html?id=81089'; I='//i1.hdslb.com/bfs/openplatform/202401/6yeUmiu11705646392215.jpeg' },
    @{ Row=13; B='2024-02-16'; C='上高·星语动漫嘉年华'; D='镜山大道2号 迎宾馆大酒店'; E='2024.02.16 09:30-02.16 17:00'; F=121; G=40; H='https://show.bilibili.com/platform/detail.html?id=80844'; I='//i1.hdslb.com/bfs/openplatform/202401/QCJN9j8h1705306410081.png' },
    @{ Row=14; B='2024-02-16'; C='南昌·运动番only'; D='南龙蟠街666号 融创茂'; E='2024.02.16 10:00-02.16 17:00'; F=290; G=60; H='https://show.bilibili.com/platform/detail.html?id=80757'; I='//i2.hdslb.com/bfs/openplatform/202401/QXLfgq7f1706180123892.jpeg' },
    @{ Row=15; B='2024-02-17'; C='九江·ACD动漫游戏嘉年华02'; D='九瑞大道与重庆路交汇处西南角 九江国际会展中心'; E='2024.02.17 10:00-02.17 17:00'; F=327; G=55; H='https://show.bilibili.com/platform/detail.html?id=81055'; I='//i0.hdslb.com/bfs/openplatform/202401/7BLpSOEZ1705574359625.jpeg' },
    @{ Row=16; B='2024-02-17'; C='江西·樟树静卿国风动漫文化展览会'; D='樟树市杏佛路89号 银河国际酒店'; E='2024.02.17 09:00-02.17 17:00'; F=334; G=40; H='https://show.bilibili.com/platform/detail.html?id=80795'; I='//i2.hdslb.com/bfs/openplatform/202401/DWQnrbtu1705044465383.jpeg' },
    @{ Row=17; B='2024-02-17'; C='赣州·第一届喵喵鱼动漫游戏展'; D='105国道东100米 毅德国际会展中心'; E='2024.02.17 09:30-02.18 16:00'; F=1780; G=50; H='https://show.bilibili.com/platform/detail.html?id=78362'; I='//i0.hdslb.com/bfs/openplatform/202311/KXRHxTLL1699521247861.png' },
    @{ Row=18; B='2024-02-18'; C='万载·第七届馨缘动漫文化展'; D='康乐街道阳乐大道217号 龙凤大酒店'; E='2024.02.18 09:30-02.18 17:00'; F=74; G=40; H='https://show.bilibili.com/platform/detail.html?id=80971'; I='//i1.hdslb.com/bfs/openplatform/202401/6ZDl6Oou1705487204077.png' },
    @{ Row=19; B='2024-02-18'; C='奉新·COP动漫游戏嘉年华1.0'; D='应星北大道482号 金勺宴大酒店'; E='2024.02.18 09:00-02.18 17:00'; F=109; G=30; H='https://show.bilibili.com/platform/detail.html?id=78259'; I='//i0.hdslb.com/bfs/openplatform/202311/yqw3kAkh1699597195072.jpeg' },
    @{ Row=20; B='2024-02-20'; C='江西·高安首届静卿国风动漫文化展览会'; D='华林中路606号 华鼎国际大酒店'; E='2024.02.20 09:00-02.20 17:00'; F=183; G=40; H='https://show.bilibili.com/platform/detail.html?id=80785'; I='//i0.hdslb.com/bfs/openplatform/202401/kcU6CEz91705040408216.jpeg' },
    @{ Row=21; B='2024-02-23'; C='上饶·囧喵喵次元国风动漫游戏展'; D='春江北大道19号 博悦宴会艺术中心'; E='2024.02.23 09:00-02.23 17:00'; F=705; G=65; H='https://show.bilibili.com/platform/detail.html?id=80240'; I='//i0.hdslb.com/bfs/openplatform/202312/Qwh83wl31703836740097.jpeg' },
    @{ Row=22; B='2024-02-23'; C='南昌·国乙only·突破次元计划（取消）'; D='高处见美好生活公园 百家喜宴高新店'; E='2024.02.23 10:00-02.23 21:00'; F=304; G='不可售'; H='https://show.bilibili.com/platform/detail.html?id=80413'; I='//i0.hdslb.com/bfs/openplatform/202401/XvmB77wb1704252353395.jpeg' },
    @{ Row=23; B='2024-02-24'; C='南昌·Cookie动漫嘉年华-赵路专场票'; D='九龙大道1177号 南昌绿地国际博览中心'; E='2024.02.24 11:00-02.24 17:00'; F=349; G='已售罄'; H='https://show.bilibili.com/platform/detail.html?id=81769'; I='//i2.hdslb.com/bfs/openplatform/202402/DhCi2kWe1707123386859.png' },
    @{ Row=24; B='2024-02-24'; C='南昌·第一届Cookie动漫嘉年华'; D='九龙大道1177号 南昌绿地国际博览中心'; E='2024.02.24 09:00-02.24 17:00'; F=4315; G=65; H='https://show.bilibili.com/platform/detail.html?id=81033'; I='//i1.hdslb.com/bfs/openplatform/202401/P994oBkz1705562167665.png' },
    @{ Row=25; B='2024-02-24'; C='宜春·融荟城难忘今宵汉文化节'; D='宜阳大道239号 宜春融荟城'; E='2024.02.24 14:00-02.24 18:00'; F=16; G=10; H='https://show.bilibili.com/platform/detail.html?id=81690'; I='//i0.hdslb.com/bfs/openplatform/202402/ldtkc9Sp1706865634128.jpeg' },
    @{ Row=26; B='2024-02-24'; C='景德镇·陶溪川×次元文化元宵游园会（ 免费活动）'; D='新厂西路315号 陶溪川发布大厅'; E='2024.02.24 10:00-02.25 18:00'; F=298; G=30; H='https://show.bilibili.com/platform/detail.html?id=81207'; I='//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png' },
    @{ Row=27; B='2024-03-02'; C='南昌·meeting动漫游戏嘉年华'; D='南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'; E='2024.03.02 09:00-03.03 17:00'; F=1140; G=60; H='https://show.bilibili.com/platform/detail.html?id=79555'; I='//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg' },
    @{ Row=28; B='2024-03-09'; C='景德镇·江报国风动漫展 '; D='迎宾大道与寺山路交叉口东200米 陶博城'; E='2024.03.09 09:00-03.10 17:00'; F=500; G=45; H='https://show.bilibili.com/platform/detail.html?id=81362'; I='//i0.hdslb.com/bfs/openplatform/202401/ae5G3ouV1706092057911.jpeg' },
    @{ Row=29; B='2024-03-16'; C='景德镇·原神X崩铁X崩坏动漫展only'; D='陶阳南路188号 晨枫臻品酒店'; E='2024.03.16 10:00-03.16 17:00'; F=46; G=55; H='https://show.bilibili.com/platform/detail.html?id=80920'; I='//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png' },
    @{ Row=30; B='2024-03-16'; C='江西·ShiningStaR动漫游戏文化节5th'; D='高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆'; E='2024.03.16 09:30-03.17 17:00'; F=669; G=60; H='https://show.bilibili.com/platform/detail.html?id=81792'; I='//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg' },
    @{ Row=31; B='2024-03-23'; C='上饶·原×铁×崩only'; D='五三东大道42号 回禾酒店'; E='2024.03.23 10:00-03.23 17:00'; F=24; G=60; H='https://show.bilibili.com/platform/detail.html?id=81103'; I='//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg' },
    @{ Row=32; B='2024-03-23'; C='南昌·AP动漫游戏嘉年华'; D='八一桥街道青山南路118号 蓝海会展中心'; E='2024.03.23 09:00-03.24 17:00'; F=334; G=60; H='https://show.bilibili.com/platform/detail.html?id=81232'; I='//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg' },
    @{ Row=33; B='2024-03-23'; C='南昌·原X穹X崩only'; D='丰和北大道299号 新吉花园酒店'; E='2024.03.23 10:00-03.23 17:00'; F=51; G=65; H='https://show.bilibili.com/platform/detail.html?id=80807'; I='//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg' },
    @{ Row=34; B='2024-03-30'; C='南昌·CM01动漫游戏博览会'; D='怀玉山大道1315号 南昌绿地国际博览中心'; E='2024.03.30 10:00-03.31 17:00'; F=172; G=55; H='https://show.bilibili.com/platform/detail.html?id=81691'; I='//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png' },
    @{ Row=35; B='2024-03-30'; C='鹰潭·原×铁×崩only'; D='南站路24号 回禾酒店(鹰潭火车站店)'; E='2024.03.30 10:00-03.30 17:00'; F=15; G=60; H='https://show.bilibili.com/platform/detail.html?id=81097'; I='//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg' }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column B holds date-like text (e.g. "2024-02-12"); force it to stay as text
    # instead of being auto-converted into a date serial number.
    $ws.Range("B2:B35").NumberFormat = "@"

    foreach ($entry in $rowsData) {
        $r = $entry.Row
        $ws.Cells.Item($r, 2).Value = $entry.B
        $ws.Cells.Item($r, 3).Value = $entry.C
        $ws.Cells.Item($r, 4).Value = $entry.D
        $ws.Cells.Item($r, 5).Value = $entry.E
        $ws.Cells.Item($r, 6).Value = $entry.F
        $ws.Cells.Item($r, 7).Value = $entry.G
        $ws.Cells.Item($r, 8).Value = $entry.H
        $ws.Cells.Item($r, 9).Value = $entry.I
    }

    # Remove the now-obsolete last row (old row 36), shrinking the used range to A1:I35
    $ws.Rows.Item(36).Delete()
}
